$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games changed the underlying stats).
# Column F (Win) is untouched; B, C, D, E are recomputed and G = B + C + D + E.
$data = @{
    2  = @{ B = 0.6606524410359556;  C = 1.655778082260271;  D = 3.537761648806719;   E = 10.19245300693656;   G = 16.0466451790395  }
    3  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 3.537761648806719;   E = 0.4942365360607697;  G = 8.974608811992548 }
    4  = @{ B = 1.455362044514542;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697;  G = 3.754798637575387 }
    5  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 6.189590430959694 }
    6  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 3.537761648806719;   E = 0.4942365360607697;  G = 8.974608811992548 }
    7  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 22.3905356188092;    E = 10.19245300693656;   G = 37.52559925287081 }
    8  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 22.3905356188092;    E = 0.4942365360607697;  G = 27.82738278199502 }
    9  = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 22.3905356188092;    E = 0.4942365360607697;  G = 27.82738278199502 }
    10 = @{ B = 1.455362044514542;   C = 1.655778082260271;  D = 22.3905356188092;    E = 0.4942365360607697;  G = 25.99591228164478 }
    11 = @{ B = 3.286832544864788;   C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 6.189590430959694 }
    12 = @{ B = 1.455362044514542;   C = 1.655778082260271;  D = 22.3905356188092;    E = 0.4942365360607697;  G = 25.99591228164478 }
    13 = @{ B = 1.455362044514542;   C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 4.358119930609447 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
